# AutoCommit_26 апреля 2024 г. 10:01:10_SibNout2023
# Mark additional answers as correct (green "5") for several students and
# let the row totals (column J, shared formula SUM(C:I)) recalculate.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A "green" (correct-answer) formatted cell we can copy the cell style from,
# so the newly-graded cells reuse the exact same cellXf/fill that the sheet
# already uses for "5" answers, instead of Excel inventing a new one.
$greenTemplate = $ws.Range("G10")

# Row 10 (student #7): G10 was blank, already green-styled -> mark as 5.
$ws.Range("G10").Value = 5

# Row 21 (student #18): E21 was a red "0" -> becomes a green "5".
$greenTemplate.Copy()
$ws.Range("E21").PasteSpecial(-4122)
$ws.Range("E21").Value = 5

# Row 29 (student #26): F29 was a red "0" -> becomes a green "5";
# G29 was blank, already green-styled -> mark as 5.
$greenTemplate.Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("F29").Value = 5
$ws.Range("G29").Value = 5

# Row 30 (student #27): C30,D30,E30,F30 were red "0" -> become green "5";
# G30 was blank, already green-styled -> mark as 5.
$greenTemplate.Copy()
$ws.Range("C30:F30").PasteSpecial(-4122)
$ws.Range("C30").Value = 5
$ws.Range("D30").Value = 5
$ws.Range("E30").Value = 5
$ws.Range("F30").Value = 5
$ws.Range("G30").Value = 5

$excel.CutCopyMode = $false

# Move the active selection to E21, matching the author's cursor position
# after re-grading that cell.
$ws.Activate()
$ws.Range("E21").Select()
